$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared string previously used as "NAO" must now read "Celular".
# Since C2 currently points at that string ("NAO") but should become "PASS",
# first repoint C2 to "PASS", then it is safe to repurpose the old text to "Celular".
$ws.Range("C2").Value = "PASS"

# Add the new row of data using the (now free) "Celular" text.
$ws.Range("A4").Value = "Celular"

# Leave the active selection on the newly added cell.
$ws.Range("A4").Select()
